$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'lunges pad'
$ws.Range('A2').Value = 'knee protection pad'
$ws.Range('A3').Value = 'boys basketball pants'
$ws.Range('A4').Value = 'working knee pads for men'
$ws.Range('A5').Value = 'knee yoga pants'
$ws.Range('A6').Value = 'padded knee sleeve for sliding'
$ws.Range('A7').Value = 'men tight pants'
$ws.Range('A8').Value = 'protective compression wear'
$ws.Range('A9').Value = 'spandex compression shorts men'
$ws.Range('A10').Value = 'softball mens pants'
$ws.Range('A11').Value = 'compression capri'
$ws.Range('A12').Value = 'softball sliding pants youth girls'
$ws.Range('A13').Value = 'baseball pants black mens'
$ws.Range('A14').Value = 'wrestling knee pads pair'
$ws.Range('A15').Value = 'compressions tights for men'
$ws.Range('A16').Value = 'men sports leggings'
$ws.Range('A17').Value = 'compression shorts boys padded'
$ws.Range('A18').Value = 'basketball padded knee sleeve'
$ws.Range('A19').Value = 'yoga pad thick'
$ws.Range('A20').Value = 'mens leggings tall'
$ws.Range('A21').Value = 'knee protector construction'
$ws.Range('A22').Value = 'basketball leggings for girls'
$ws.Range('A23').Value = 'lacrosse compression shorts'
$ws.Range('A24').Value = 'boys sports tights leggings'
$ws.Range('A25').Value = 'mens sweat pads'
$ws.Range('A26').Value = 'mens construction knee pads'
$ws.Range('A27').Value = 'yoga pants men'
$ws.Range('A28').Value = 'bursitis knee'
$ws.Range('A29').Value = 'thread protector paintball'
$ws.Range('A30').Value = 'tight gym pants men'
$ws.Range('A31').Value = 'adult bee tights'
$ws.Range('A32').Value = 'basketball shorts in bulk'
$ws.Range('A33').Value = 'rash guard men leggings'
$ws.Range('A34').Value = 'youth sliding short'
$ws.Range('A35').Value = 'running knee compression'
$ws.Range('A36').Value = 'knee protection for running'
$ws.Range('A37').Value = 'girls basketball leggings'
$ws.Range('A38').Value = 'youth xxl baseball pants'
$ws.Range('A39').Value = 'boys softball pants'
$ws.Range('A40').Value = 'tight for boys'
$ws.Range('A41').Value = 'big and tall leggings men'
$ws.Range('A42').Value = 'knee protect'
$ws.Range('A43').Value = 'knee sleeve basketball men'
$ws.Range('A44').Value = 'protective baseball'
$ws.Range('A45').Value = 'compression pants for recovery'
$ws.Range('A46').Value = 'girls basketball knee guards'
$ws.Range('A47').Value = 'basketball knee pads for women'
$ws.Range('A48').Value = 'gel wrestling knee pads'
$ws.Range('A49').Value = 'mens hockey pads'
$ws.Range('A50').Value = 'long knee pads volleyball'
$ws.Range('A51').Value = 'mens compression pants 3/4 length'
$ws.Range('A52').Value = 'knee pads motorcycle'
$ws.Range('A53').Value = 'sliding shorts youth girls'
$ws.Range('A54').Value = 'boys athletic leggings youth'
$ws.Range('A55').Value = 'knee pads for basketball women'
$ws.Range('A56').Value = 'big mens compression pants'
$ws.Range('A57').Value = 'soccer sliding pants'
$ws.Range('A58').Value = 'good thread mens pants'
$ws.Range('A59').Value = 'xl knee pads for men'
$ws.Range('A60').Value = 'capri pouches adults'
$ws.Range('A61').Value = 'knee sleeves hex'
$ws.Range('A62').Value = 'knee compression sleeve with padding'
$ws.Range('A63').Value = 'baleaf mens pants'
$ws.Range('A64').Value = 'leggings tight'
$ws.Range('A65').Value = 'compression tights mens'
$ws.Range('A66').Value = 'compression knee sleeve basketball'
$ws.Range('A67').Value = 'volleyball knee pads gel'
$ws.Range('A68').Value = 'knee pads for work xxl'
$ws.Range('A69').Value = 'knee pad for exercise'
$ws.Range('A70').Value = 'running tights youth boys'
$ws.Range('A71').Value = 'boy capri pants'
$ws.Range('A72').Value = 'youth hockey girdle'
$ws.Range('A73').Value = 'youth compression leggings boys'
$ws.Range('A74').Value = 'sliding compression shorts'
$ws.Range('A75').Value = 'baseball leg protection'
$ws.Range('A76').Value = 'knee sleeves for basketball youth'
$ws.Range('A77').Value = 'youth paintball pants'
$ws.Range('A78').Value = 'youth boy tights'
$ws.Range('A79').Value = 'softball catcher pants'
$ws.Range('A80').Value = 'knee pads volleyball mens'
$ws.Range('A81').Value = 'mens running compression pants'
$ws.Range('A82').Value = 'football tights youth boys'
$ws.Range('A83').Value = 'athletic tights youth boys'
$ws.Range('A84').Value = 'padded calf sleeve'
$ws.Range('A85').Value = 'mens softball pants'
$ws.Range('A86').Value = 'mens baseball shorts'
$ws.Range('A87').Value = 'knee pad for soccer'
$ws.Range('A88').Value = 'compression leggings boys'
$ws.Range('A89').Value = 'capri for men'
$ws.Range('A90').Value = 'mens workout pants leggings'
$ws.Range('A91').Value = 'knee pad sport'
$ws.Range('A92').Value = 'tights pants boys'
$ws.Range('A93').Value = 'softball sliding shorts women'
$ws.Range('A94').Value = 'padded sliding shorts'
$ws.Range('A95').Value = 'arthritis hope knee sleeve'
$ws.Range('A96').Value = 'padded knee sleeve basketball'
$ws.Range('A97').Value = 'snowboarding mens pants'
$ws.Range('A98').Value = 'knee compression cold'
$ws.Range('A99').Value = 'knee pads squats'
$ws.Range('A100').Value = 'men work pants with knee pads'
